$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously empty grade cells for row 16 (student #13)
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 5

# Match the formatting already used in the neighbouring grade columns
$ws.Range("I17").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("J15").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active cell / selection to match the new state
$ws.Range("J16").Select()
